$d = $word.ActiveDocument

# Common namespace declaration reused for every InsertXML fragment.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Paragraph 1 ("Show loader with this line ...") : highlight yellow -> green
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$xml1 = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Show loader with </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">this line “the answers lie within you” </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>for the first time whenever a new user come to the website.</w:t></w:r></w:p>
"@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Paragraph 2 ("Add fade in and fade out effect on loader ...") : add green
# highlight on the paragraph mark and every run (no text/structure change).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$xml2 = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Add fade in and fade out effect on loader</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> for </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>2-</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">3 </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>)</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>seconds</w:t></w:r><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>.</w:t></w:r></w:p>
"@
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Paragraph 3 ("Change login page color.") : add green highlight
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$xml3 = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Change login page color.</w:t></w:r></w:p>
"@
$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Paragraph 4 ("Make some UI for onboarding page.") : merge the two runs
# into one, drop the _GoBack bookmark (it moves to paragraph 5 below), and
# add a yellow highlight.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$xml4 = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Make some UI for onboarding page.</w:t></w:r></w:p>
"@
$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------------
# Paragraph 5 ("Change hover login button color.") : split the single run
# into two, insert the _GoBack bookmark between them, and add a green
# highlight.
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$xml5 = @"
<w:p $wns><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Change hover login button </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>color.</w:t></w:r></w:p>
"@
$p5.Range.InsertXML($xml5)

Write-Host "edits applied"
